# Add cantrals by cantons
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Move existing data rows 3,4,5 up to 2,3,4 (values only change row index) ---
# Row values currently at rows 3-5 (A:K) need to become rows 2-4.
# Easiest: capture current values for rows 3,4,5 then delete old rows 1 & 2
# and rewrite row 1 as new header, rewriting data rows 2-4 from the captured data.

$names = @("La Goule", "Bellerive", "Bassecourt")
$data = @(
    @(1, 510100, $names[0], 1894, 1958, 22, 5.6, 5, 14, 12, 26),
    @(2, 109915, $names[1], 1905, 2002, 9.8, 0.52, 0.46, 1.34, 0.96, 2.3),
    @(3, 109900, $names[2], 1920, 2001, 3, 1.03, 0.86, 2.6, 1.9, 4.5)
)

# Clear the whole used range first so we fully control the final layout
$ws.Cells.Clear()

# --- Header row (row 1) ---
$headers = @("idx", "idx2", "Name", "Date Start", "Date End", "(m3/s)", "(MW1)", "(MW2)", "(GWh) Winter", "(GWh) Summer", "(GWh) Year")
for ($c = 1; $c -le $headers.Length; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 1]
    if ($c -ge 6) {
        $cell.Style = "Normal"
    }
}

# --- Data rows (rows 2..4) ---
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c).Value = $row[$c - 1]
    }
}

$ws.Range("A2:K2").Select()
